# Automatische test-sync: 2025-06-22 19:00:50
#
# Adds a new incoming mail log entry (row 30) to the "Logs" sheet and
# refreshes the "Dashboard" category summary table to reflect the
# updated counts/ordering.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append the new row (row 30)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A30").Value = "Zakelijke samenwerking"
$logs.Range("B30").Value = "mailmind.test@zohomail.eu"
$logs.Range("C30").Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D30").Value = "Samenwerking / Partnerverzoek"
$logs.Range("F30").Value = "2025-06-22 19:00:14"
$logs.Range("G30").Value = "Nee"

# Extend the conditional formatting ranges so they cover the new row,
# mirroring the dimension growing from A1:G29 to A1:G30.
$catRules = $logs.Range("D2:D29").FormatConditions
for ($i = 1; $i -le $catRules.Count; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D30"))
}

$answeredRules = $logs.Range("G2:G29").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G30"))
}

# ---------------------------------------------------------------------
# 2) Dashboard sheet: re-sync the category/count summary table so the
#    newly added "Samenwerking / Partnerverzoek" entry is reflected and
#    the rows are re-ordered to match the refreshed counts.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Retour / Terugbetaling"
$dash.Range("B4").Value = 3

$dash.Range("A5").Value = "Afmelding / Nieuwsbrief"
$dash.Range("B5").Value = 3

$dash.Range("A6").Value = "Samenwerking / Partnerverzoek"
$dash.Range("B6").Value = 3

$dash.Range("A9").Value = "Offerte / Prijsaanvraag"
$dash.Range("B9").Value = 2
